{"js": "const replacements = [\n  [\"2025-02-24 Monday\", \"2025-02-25 Tuesday\"],\n  [\"604\u00f77=\", \"686\u00f79=\"],\n  [\"357\u00f78=\", \"576\u00f74=\"],\n  [\"407\u00f73=\", \"253\u00f77=\"],\n  [\"149\u00f75=\", \"497\u00f72=\"],\n  [\"231\u00f78=\", \"470\u00f78=\"],\n  [\"624\u00f78=\", \"826\u00f76=\"],\n  [\"538\u00f75=\", \"161\u00f73=\"],\n  [\"675\u00f73=\", \"131\u00f72=\"],\n  [\"810\u00f78=\", \"788\u00f79=\"],\n  [\"714\u00f78=\", \"346\u00f77=\"],\n  [\"721\u00f78=\", \"598\u00f73=\"],\n  [\"838\u00f73=\", \"818\u00f74=\"],\n  [\"884\u00f78=\", \"224\u00f78=\"],\n  [\"906\u00f77=\", \"962\u00f76=\"],\n  [\"746\u00f74=\", \"480\u00f75=\"],\n  [\"390\u00f79=\", \"239\u00f79=\"],\n  [\"467\u00f77=\", \"836\u00f79=\"],\n  [\"407\u00f72=\", \"731\u00f73=\"],\n  [\"112\u00f75=\", \"782\u00f73=\"],\n  [\"130\u00f77=\", \"446\u00f73=\"],\n  [\"929\u00f74=\", \"779\u00f73=\"],\n  [\"949\u00f79=\", \"968\u00f78=\"],\n  [\"945\u00f73=\", \"207\u00f78=\"],\n  [\"988\u00f77=\", \"680\u00f73=\"],\n  [\"235\u00f78=\", \"629\u00f75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-24 Monday\", \"2025-02-25 Tuesday\"),\n    @(\"604\u00f77=\", \"686\u00f79=\"),\n    @(\"357\u00f78=\", \"576\u00f74=\"),\n    @(\"407\u00f73=\", \"253\u00f77=\"),\n    @(\"149\u00f75=\", \"497\u00f72=\"),\n    @(\"231\u00f78=\", \"470\u00f78=\"),\n    @(\"624\u00f78=\", \"826\u00f76=\"),\n    @(\"538\u00f75=\", \"161\u00f73=\"),\n    @(\"675\u00f73=\", \"131\u00f72=\"),\n    @(\"810\u00f78=\", \"788\u00f79=\"),\n    @(\"714\u00f78=\", \"346\u00f77=\"),\n    @(\"721\u00f78=\", \"598\u00f73=\"),\n    @(\"838\u00f73=\", \"818\u00f74=\"),\n    @(\"884\u00f78=\", \"224\u00f78=\"),\n    @(\"906\u00f77=\", \"962\u00f76=\"),\n    @(\"746\u00f74=\", \"480\u00f75=\"),\n    @(\"390\u00f79=\", \"239\u00f79=\"),\n    @(\"467\u00f77=\", \"836\u00f79=\"),\n    @(\"407\u00f72=\", \"731\u00f73=\"),\n    @(\"112\u00f75=\", \"782\u00f73=\"),\n    @(\"130\u00f77=\", \"446\u00f73=\"),\n    @(\"929\u00f74=\", \"779\u00f73=\"),\n    @(\"949\u00f79=\", \"968\u00f78=\"),\n    @(\"945\u00f73=\", \"207\u00f78=\"),\n    @(\"988\u00f77=\", \"680\u00f73=\"),\n    @(\"235\u00f78=\", \"629\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2) | Out-Null\n}\n"}
